$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Sweden"
$ws.Range("B11").Value = 129.7
$ws.Range("C11").Value = 131.6
$ws.Range("D11").Value = 120
$ws.Range("E11").Value = 127.1
$ws.Range("F11").Value = 134
$ws.Range("G11").Value = 173.3
$ws.Range("H11").Value = 181.5
$ws.Range("I11").Value = 248.9
$ws.Range("J11").Value = 178.4
$ws.Range("K11").Value = 140.7
$ws.Range("L11").Value = 157.1
$ws.Range("M11").Value = 150.7
$ws.Range("N11").Value = 172.7
